# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 00:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 529357
$ws.Range("C4").Value = 26481
$ws.Range("E4").Value = 479448
$ws.Range("G4").Value = 1720
$ws.Range("H4").Value = 20467

# Alemania (row 8)
$ws.Range("B8").Value = 124908
$ws.Range("C8").Value = 2737
$ws.Range("E8").Value = 64772

# China (row 9)
$ws.Range("C9").Value = 0
$ws.Range("G9").Value = 0

# Reorder: Bahamas now sorts before Islas Caimanes / Macao (rows 149-151)
$ws.Range("A149").Value = "Bahamas"
$ws.Range("B149").Value = 46
$ws.Range("C149").Value = 4
$ws.Range("D149").Value = 5
$ws.Range("E149").Value = 33
$ws.Range("F149").Value = 1
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 8

$ws.Range("A150").Value = "Islas Caimanes"
$ws.Range("B150").Value = 45
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 6
$ws.Range("E150").Value = 38
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 1

$ws.Range("A151").Value = "Macao"
$ws.Range("B151").Value = 45
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 10
$ws.Range("E151").Value = 35
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0

# Row 152 (Guyana) is unaffected

# Reorder: Zimbabue now sorts before Curazao / Botsuana / Belice (rows 181-184)
$ws.Range("A181").Value = "Zimbabue"
$ws.Range("B181").Value = 14
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 11
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 3

$ws.Range("A182").Value = "Curazao"
$ws.Range("B182").Value = 14
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 7
$ws.Range("E182").Value = 6
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 1

$ws.Range("A183").Value = "Botsuana"
$ws.Range("B183").Value = 13
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 12
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 1

$ws.Range("A184").Value = "Belice"
$ws.Range("B184").Value = 13
$ws.Range("C184").Value = 3
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 11
$ws.Range("F184").Value = 1
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 2
